# The deck contained two back-to-back duplicate "Agenda (This should be
# the PPT flow)" slides (presentation order 11 and 12, slide IDs 267 and
# 269 respectively). This edit removes the second, redundant copy
# (slide id 269 / the slide at presentation position 12), which also
# drops its associated notes page. All following slides shift up by one
# position; their slide IDs and content are otherwise untouched.

$p = $ppt.ActivePresentation

$agendaIndexes = @()
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $titleText = ""
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $titleText = $shp.TextFrame.TextRange.Text
            break
        }
    }
    if ($titleText -like "Agenda (This should be the PPT flow)*") {
        $agendaIndexes += $i
    }
}

if ($agendaIndexes.Count -ge 2) {
    # Remove the second (duplicate) occurrence.
    $targetIndex = $agendaIndexes[1]
} elseif ($agendaIndexes.Count -eq 1) {
    $targetIndex = $agendaIndexes[0]
} else {
    $targetIndex = 12
}

$p.Slides.Item($targetIndex).Delete()
